$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block (rows 13-15): secondary "tare" scratch computation ---
# Written first so the new shared strings land in the same order the
# final workbook expects ("T", "Av" before "Tare"/"DHT22"/the note).
$ws.Range("A13").Value = "T"
$ws.Range("B13").Value = "Av"

$ws.Range("A14").Value = -21000
$ws.Range("B14").Value = -75826
$ws.Range("D14").Formula = "=B14-A14"
$ws.Range("F14").Formula = "=D14*-1"
$ws.Range("H14").Formula = "=D14/I7"

$ws.Range("H15").Formula = "=F14/I7"

# --- Row 18: free-form note left by the author ---
$ws.Range("G18").Value = "{-21000,-23208.92,20,0},     // J03 évolution valeurs en négatif. tester sur bornier"

# --- Row 1: C1 becomes the "Tare" label (was a literal 0) ---
$ws.Range("C1").Value = "Tare"

# --- Row 10 ---
$ws.Range("A10").Value = 6
$ws.Range("J10").Value = "DHT22"

# --- Row 5: G5's shared formula gets "detached" into its own literal formula ---
$ws.Range("G5").Formula = "=E5-C5"
$ws.Range("G5").Formula = "=+E5-C5"

# --- Row 6: new shared-formula anchor (G6:G8) + ratio column ---
$ws.Range("G6:G8").Formula = "=+E6-C6"
$ws.Range("I6").Formula = "=G6/poids"

# --- Row 7 ---
$ws.Range("A7").Value = 3
$ws.Range("C7").Value = -21800
$ws.Range("C7").Interior.Color = 65535
$ws.Range("E7").Value = -75900
$ws.Range("E7").Interior.Color = 65535
$ws.Range("I7").Formula = "=G7/poids"

# --- Row 8 ---
$ws.Range("A8").Value = 8
$ws.Range("C8").Value = -35751
$ws.Range("C8").Interior.Color = 65535
$ws.Range("E8").Value = -88863
$ws.Range("E8").Interior.Color = 65535
$ws.Range("I8").Formula = "=G8/poids"

# --- Row 9 ---
$ws.Range("A9").Value = 9
$ws.Range("C9").Value = -28026
$ws.Range("C9").Interior.Color = 65535
$ws.Range("E9").Value = -81617
$ws.Range("E9").Interior.Color = 65535
$ws.Range("G9").Formula = "=+E9-C9"
$ws.Range("I9").Formula = "=G9/poids"

# --- Column width for I (raw OOXML width 14 == COM ColumnWidth 14 - 5/6) ---
$ws.Columns.Item(9).ColumnWidth = 13.166666666666666

# --- Selection state ---
$ws.Range("I9").Select()
